$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "SlowAntibodies" passive skill as row 23, mirroring the other
# zero-cost passive skills (PassiveSkill type, 0 Mana, 0 CoolDown).
# Inserting the row (rather than just writing into a blank one) makes Excel
# inherit the same cell formatting ("Bom"/Good style) used by the rows above.
$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "SlowAntibodies"
$ws.Range("C23").Value = "PassiveSkill"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0

# Mirror the diff's updated selection (the new last row, whole-row selected)
$ws.Rows.Item(23).Select()
